$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 34743.668
$ws.Range("I69").Value = 5000
$ws.Range("J69").Value = 37447.637
$ws.Range("K69").Value = 15000
$ws.Range("L69").Value = 112342.911
$ws.Range("M69").Value = -14126
$ws.Range("N69").Value = -114090.911
$ws.Range("H72").Value = 34743.668
$ws.Range("I72").Value = 5000
$ws.Range("J72").Value = 37447.637
$ws.Range("K72").Value = 45000
$ws.Range("L72").Value = 337028.733
$ws.Range("M72").Value = -40632
$ws.Range("N72").Value = -345764.733
$ws.Range("H74").Value = 9160.200000000001
$ws.Range("I74").Value = 6800.4287
$ws.Range("K74").Value = 6800.4287
$ws.Range("M74").Value = -5864.4287
$ws.Range("H77").Value = 9160.200000000001
$ws.Range("I77").Value = 6800.4287
$ws.Range("K77").Value = 34002.14350000001
$ws.Range("M77").Value = -29322.14350000001
$ws.Range("H80").Value = 1857.2195
$ws.Range("J80").Value = 2786.0435
$ws.Range("L80").Value = 8358.130500000001
$ws.Range("N80").Value = -10354.1305
$ws.Range("H83").Value = 1857.2195
$ws.Range("J83").Value = 2786.0435
$ws.Range("L83").Value = 25074.3915
$ws.Range("N83").Value = -35058.3915
$ws.Range("H87").Value = 22081.791
$ws.Range("J87").Value = 22081.791
$ws.Range("L87").Value = 22081.791
$ws.Range("N87").Value = -24577.791
$ws.Range("H90").Value = 22081.791
$ws.Range("J90").Value = 22081.791
$ws.Range("L90").Value = 66245.37300000001
$ws.Range("N90").Value = -78725.37300000001
$ws.Range("H98").Value = 856.25
$ws.Range("I98").Value = 880.3333
$ws.Range("K98").Value = 880.3333
$ws.Range("M98").Value = 617.6667
$ws.Range("H122").Value = 856.25
$ws.Range("I122").Value = 880.3333
$ws.Range("K122").Value = 2640.9999
$ws.Range("M122").Value = -190.9998999999998
$ws.Range("H137").Value = 3841.0833
$ws.Range("I137").Value = 2834.2942
$ws.Range("J137").Value = 6286.143
$ws.Range("K137").Value = 8502.882599999999
$ws.Range("L137").Value = 18858.429
$ws.Range("M137").Value = -5952.882599999999
$ws.Range("N137").Value = -23958.429
$ws.Range("H138").Value = 5199.421
$ws.Range("J138").Value = 5456.75
$ws.Range("L138").Value = 16370.25
$ws.Range("N138").Value = -26650.25
$ws.Range("H141").Value = 6227.231
$ws.Range("I141").Value = 6329.5
$ws.Range("K141").Value = 18988.5
$ws.Range("M141").Value = -13808.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 14974.916
$ws.Range("I32").Value = 11400.092
$ws.Range("J32").Value = 29274.21
$ws.Range("K32").Value = 11400.092
$ws.Range("L32").Value = 29274.21
$ws.Range("M32").Value = -11113.092
$ws.Range("N32").Value = -29848.21
$ws.Range("H135").Value = 68173.42999999999
$ws.Range("J135").Value = 68173.42999999999
$ws.Range("L135").Value = 68173.42999999999
$ws.Range("N135").Value = -78313.42999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 92806
$ws.Range("J20").Value = 252075.25
$ws.Range("L20").Value = 252075.25
$ws.Range("N20").Value = -252569.25
$ws.Range("H81").Value = 51390
$ws.Range("J81").Value = 51390
$ws.Range("L81").Value = 51390
$ws.Range("N81").Value = -53512
$ws.Range("H84").Value = 51390
$ws.Range("J84").Value = 51390
$ws.Range("L84").Value = 154170
$ws.Range("N84").Value = -164778

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3631.5103
$ws.Range("I31").Value = 2685.9355
$ws.Range("J31").Value = 5260
$ws.Range("K31").Value = 2685.9355
$ws.Range("L31").Value = 5260
$ws.Range("M31").Value = -2390.9355
$ws.Range("N31").Value = -5850
$ws.Range("H34").Value = 3631.5103
$ws.Range("I34").Value = 2685.9355
$ws.Range("J34").Value = 5260
$ws.Range("K34").Value = 2685.9355
$ws.Range("L34").Value = 5260
$ws.Range("M34").Value = -2483.9355
$ws.Range("N34").Value = -5664
$ws.Range("H99").Value = 14030943
$ws.Range("I99").Value = 2223917.2
$ws.Range("K99").Value = 2223917.2
$ws.Range("M99").Value = -2222419.2
$ws.Range("H126").Value = 14030943
$ws.Range("I126").Value = 2223917.2
$ws.Range("K126").Value = 6671751.600000001
$ws.Range("M126").Value = -6669281.600000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H38").Value = 25997
$ws.Range("J38").Value = 25997
$ws.Range("L38").Value = 25997
$ws.Range("N38").Value = -26923
$ws.Range("H57").Value = 18166.334
$ws.Range("H70").Value = 127719.336
$ws.Range("I70").Value = 225179.8
$ws.Range("J70").Value = 5893.75
$ws.Range("K70").Value = 225179.8
$ws.Range("L70").Value = 5893.75
$ws.Range("M70").Value = -224909.8
$ws.Range("N70").Value = -6433.75
$ws.Range("H73").Value = 127719.336
$ws.Range("I73").Value = 225179.8
$ws.Range("J73").Value = 5893.75
$ws.Range("K73").Value = 225179.8
$ws.Range("L73").Value = 5893.75
$ws.Range("M73").Value = -224243.8
$ws.Range("N73").Value = -7765.75
$ws.Range("H80").Value = 33447678
$ws.Range("I80").Value = 188759.83
$ws.Range("K80").Value = 188759.83
$ws.Range("M80").Value = -187761.83
$ws.Range("H83").Value = 33447678
$ws.Range("I83").Value = 188759.83
$ws.Range("K83").Value = 943799.1499999999
$ws.Range("M83").Value = -938807.1499999999
$ws.Range("H126").Value = 3924.0588
$ws.Range("I126").Value = 2289.6
$ws.Range("K126").Value = 6868.799999999999
$ws.Range("M126").Value = -4398.799999999999
$ws.Range("H135").Value = 61666.5
$ws.Range("J135").Value = 61666.5
$ws.Range("L135").Value = 61666.5
$ws.Range("N135").Value = -71806.5
$ws.Range("H141").Value = 96106.75
$ws.Range("J141").Value = 96106.75
$ws.Range("L141").Value = 96106.75
$ws.Range("N141").Value = -106466.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5514.091
$ws.Range("I7").Value = 2925
$ws.Range("K7").Value = 2925
$ws.Range("M7").Value = -2813
$ws.Range("H22").Value = 1410.1
$ws.Range("I22").Value = 1333.3334
$ws.Range("K22").Value = 1333.3334
$ws.Range("M22").Value = -1038.3334
$ws.Range("H27").Value = 1410.1
$ws.Range("I27").Value = 1333.3334
$ws.Range("K27").Value = 1333.3334
$ws.Range("M27").Value = -1226.3334
$ws.Range("H40").Value = 14976.75
$ws.Range("I40").Value = 24702.2
$ws.Range("K40").Value = 24702.2
$ws.Range("M40").Value = -24566.2
$ws.Range("I61").Value = 2030.6
$ws.Range("J61").Value = 8000
$ws.Range("K61").Value = 2030.6
$ws.Range("L61").Value = 8000
$ws.Range("M61").Value = -1828.6
$ws.Range("N61").Value = -8404
$ws.Range("H68").Value = 9541.166999999999
$ws.Range("I68").Value = 9560
$ws.Range("J68").Value = 9527.714
$ws.Range("K68").Value = 9560
$ws.Range("L68").Value = 9527.714
$ws.Range("M68").Value = -8811
$ws.Range("N68").Value = -11025.714
$ws.Range("H71").Value = 9541.166999999999
$ws.Range("I71").Value = 9560
$ws.Range("J71").Value = 9527.714
$ws.Range("K71").Value = 47800
$ws.Range("L71").Value = 47638.57
$ws.Range("M71").Value = -44056
$ws.Range("N71").Value = -55126.57
$ws.Range("I113").Value = 2030.6
$ws.Range("J113").Value = 8000
$ws.Range("K113").Value = 2030.6
$ws.Range("L113").Value = 8000
$ws.Range("M113").Value = 139.4000000000001
$ws.Range("N113").Value = -12340
$ws.Range("H122").Value = 6942.3076
$ws.Range("I122").Value = 4781.25
$ws.Range("K122").Value = 14343.75
$ws.Range("M122").Value = -11893.75
$ws.Range("H126").Value = 5514.091
$ws.Range("I126").Value = 2925
$ws.Range("K126").Value = 8775
$ws.Range("M126").Value = -6305

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2591.3333
$ws.Range("I126").Value = 2669.6
$ws.Range("K126").Value = 8008.799999999999
$ws.Range("M126").Value = -5538.799999999999
$ws.Range("H140").Value = 47400
$ws.Range("J140").Value = 29500
$ws.Range("L140").Value = 29500
$ws.Range("N140").Value = -39860
$ws.Range("H141").Value = 99899.8
$ws.Range("J141").Value = 99899.8
$ws.Range("L141").Value = 99899.8
$ws.Range("N141").Value = -110259.8
